$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows(23).Insert()
$ws.Range("A23").Value = 43380
